$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 72, shifting existing rows 72-159 down to 73-160.
$ws.Rows("72:72").Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)

# Populate the newly inserted row 72 with its data.
$ws.Range("A72").Value = 4
$ws.Range("B72").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C72").Value = "Los Lagos"
$ws.Range("D72").Value = 44482
$ws.Range("E72").Value = 10
$ws.Range("F72").Value = 100112044
$ws.Range("G72").Value = "Perejil"
$ws.Range("H72").Value = "Sin especificar"
$ws.Range("I72").Value = "Primera"
$ws.Range("J72").Value = 30
$ws.Range("K72").Value = 4500
$ws.Range("L72").Value = 4500
$ws.Range("M72").Value = 4500
$ws.Range("N72").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O72").Value = "Región Metropolitana"
$ws.Range("P72").Value = 1500
$ws.Range("Q72").Value = 3
$ws.Range("R72").Value = "Hortaliza"
